$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to keep a plain text value (matching the source
    # workbook's inlineStr cells) instead of letting Excel auto-convert
    # numeric-looking strings (e.g. "289.83", "6.990") into floating point
    # numbers, which would both change the cell type and lose precision
    # (trailing zeros, thousand-dot formatted figures, etc).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row -> [D new value (or $null to leave unchanged), E new value (or $null)]
$updates = @(
    @(2,  "22.475.41", "  +0.46%  "),
    @(3,  "1.571.12",  "  +0.22%  "),
    @(4,  $null,        "  -0.42%  "),
    @(5,  $null,        "  -0.32%  "),
    @(6,  "289.83",    "  +0.43%  "),
    @(7,  "0.3689",    "  -1.17%  "),
    @(8,  "50.29",     "  +2.45%  "),
    @(9,  $null,        "  +0.83%  "),
    @(10, $null,        "  +2.62%  "),
    @(11, "0.07527",   "  +2.12%  "),
    @(12, $null,        "  -0.41%  "),
    @(13, "21.17",     "  +2.45%  "),
    @(14, "6.014",     "  +2.48%  "),
    @(15, "6.990",     "  +2.51%  "),
    @(16, "1.570.91",  "  +0.33%  "),
    @(17, $null,        "  +1.45%  "),
    @(18, "90.22",     "  +1.67%  "),
    @(19, "0.06757",   "  +1.19%  "),
    @(20, $null,        "  -0.44%  "),
    @(21, "6.355",     "  +3.94%  "),
    @(22, "16.32",     "  +1.74%  "),
    @(23, $null,        "  +3.17%  "),
    @(24, "22.467.36", "  +0.48%  "),
    @(25, "2.393",     "  +1.41%  "),
    @(26, "2.640",     "  +4.37%  "),
    @(27, "19.93",     "  +0.52%  "),
    @(28, "149.56",    "  +1.83%  "),
    @(29, "5.057",     "  +1.43%  "),
    @(30, "124.60",    "  +0.22%  "),
    @(31, "1.745.54",  "  +0.26%  "),
    @(32, "1.056",     "  +8.06%  "),
    @(33, "6.235",     "  +6.10%  "),
    @(34, "2.015",     "  +1.56%  "),
    @(35, "9.802",     "  +3.30%  "),
    @(36, "0.08344",   "  +0.32%  "),
    @(37, "0.02470",   "  +1.48%  "),
    @(38, "0.2293",    "  +2.53%  "),
    @(39, "1.343",     "  -3.14%  "),
    @(40, "0.06508",   $null),
    @(41, "5.414",     "  +1.26%  "),
    @(42, $null,        "  +2.56%  "),
    @(43, "0.6217",    "  +0.90%  "),
    @(46, "3.781",     "  +0.03%  "),
    @(47, "0.5847",    "  +2.00%  "),
    @(48, "2.056",     "  +1.28%  "),
    @(49, "125.67",    "  +0.27%  "),
    @(50, "1.235",     "  +1.63%  "),
    @(51, "0.07317",   "  +0.38%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]
    if ($null -ne $dVal) {
        Set-TextValue $ws.Cells.Item($row, 4) $dVal
    }
    if ($null -ne $eVal) {
        Set-TextValue $ws.Cells.Item($row, 5) $eVal
    }
}

# Rows 44 and 45 swap their coin content (EnergySwap <-> Frax) along with
# refreshed price/volume figures.
$ws.Cells.Item(44, 2).Value = "Frax"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Cells.Item(44, 4) "1.001"
Set-TextValue $ws.Cells.Item(44, 5) "  -0.44%  "

$ws.Cells.Item(45, 2).Value = "EnergySwap"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Cells.Item(45, 4) "14.05"
Set-TextValue $ws.Cells.Item(45, 5) "  +1.87%  "
